# edit.ps1 - apply the teachercom_presentation.pptx changes described by the
# commit "Added the latest changes to the powerpoints":
#
#   Slide 3 (Content Placeholder 2):
#     - "Automation through scheduling" -> "One-step student import"
#
#   Slide 4 (Content Placeholder 2):
#     - "Schools have a limited IT resources" split into three runs and the
#       stray "a " dropped: "Schools " / "have " / "limited IT resources"
#     - "Legal Issues concerning Student data and " + "parental consent"
#       (two runs) merged back into a single run.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 3 - Content Placeholder 2
# ---------------------------------------------------------------------
$slide3 = $p.Slides.Item(3)
$body3 = $slide3.Shapes.Item(2).TextFrame.TextRange

# Paragraph layout (1-based char offsets) before the edit:
#   1   : "Fast, easy to use online notification system"  (44 chars)
#   46  : "Different types of communication"               (32 chars)
#   79  : "Automation through scheduling"                  (29 chars)
#   109 : "Low-cost technology"                             (19 chars)
$body3.Characters(79, 29).Text = "One-step student import"

# ---------------------------------------------------------------------
# Slide 4 - Content Placeholder 2
# ---------------------------------------------------------------------
$slide4 = $p.Slides.Item(4)
$body4 = $slide4.Shapes.Item(2).TextFrame.TextRange

# Paragraph 1, "Schools have a limited IT resources" (chars 1-35), becomes
# three runs: "Schools " / "have " / "limited IT resources" (drops "a ").
# Edit back-to-front so earlier offsets stay valid while the text length
# changes (35 -> 33 chars overall).
$body4.Characters(14, 22).Text = "limited IT resources"
$body4.Characters(9, 5).Text = "have "
$body4.Characters(1, 8).Text = "Schools "

# Paragraph 5 used to be two runs, "Legal Issues concerning Student data and "
# + "parental consent" (originally chars 165-221, 57 chars total). The three
# edits above shortened the text before it by 2 chars ("a " was dropped), so
# paragraph 5 now starts at 163. Merge the two runs into one run with the
# same combined text.
$body4.Characters(163, 57).Text = "Legal Issues concerning Student data and parental consent"
